$wb = $excel.ActiveWorkbook

# --- 1. Insert the new "2022-Q4" sheet right after "总计", before "2021-Q3" ---
$totalSheet = $wb.Worksheets.Item(1)
$q3Sheet    = $wb.Worksheets.Item(2)

$newSheet = $wb.Worksheets.Add($q3Sheet)
$newSheet.Name = "2022-Q4"

# Header row (row 1), columns B..H
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
$col = 2
foreach ($h in $headers) {
    $newSheet.Cells.Item(1, $col).Value = $h
    $col = $col + 1
}

# Data rows
$data = @(
    @(0, "004945", "长信中证500指数增强A", "2.08", "92.79", "1.59", "0.0331", 10),
    @(1, "013881", "长信中证500指数增强C", "0.96", "92.79", "1.59", "0.0153", 10),
    @(2, "007943", "富安达中证 500 指数增强", "0.30", "78.45", "1.12", "0.0034", 8)
)

$row = 2
foreach ($d in $data) {
    $newSheet.Cells.Item($row, 1).Value = $d[0]
    $newSheet.Cells.Item($row, 2).Value = $d[1]
    $newSheet.Cells.Item($row, 3).Value = $d[2]
    $newSheet.Cells.Item($row, 4).Value = $d[3]
    $newSheet.Cells.Item($row, 5).Value = $d[4]
    $newSheet.Cells.Item($row, 6).Value = $d[5]
    $newSheet.Cells.Item($row, 7).Value = $d[6]
    $newSheet.Cells.Item($row, 8).Value = $d[7]
    $row = $row + 1
}

# Formatting to match the other quarter sheets: bold, centered header row
# and bold centered index column (A2:A4), with thin borders.
$headerRange = $newSheet.Range("B1:H1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

$idxRange = $newSheet.Range("A2:A4")
$idxRange.Font.Bold = $true
$idxRange.HorizontalAlignment = -4108
$idxRange.VerticalAlignment = -4160
$idxRange.Borders.LineStyle = 1

# --- 2. Update the "总计" sheet: insert a new row for 2022-Q4 above 2021-Q3 ---
$totalSheet.Rows.Item(2).Insert()

# The inserted row picks up stray formatting from the row above; reset the
# data columns to the plain/default style used throughout this sheet.
$totalSheet.Range("B2:D2").Style = "Normal"

# Give the new index cell (A2) the same bold/centered/bordered look as the
# other index-column cells (A3, A4) by copying their format over.
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

$totalSheet.Cells.Item(2, 1).Value = 0
$totalSheet.Cells.Item(2, 2).Value = "2022-Q4"
$totalSheet.Cells.Item(2, 3).Value = 3
$totalSheet.Cells.Item(2, 4).Value = 0.05

# Renumber the index column (A) for the rows that shifted down
$totalSheet.Cells.Item(3, 1).Value = 1
$totalSheet.Cells.Item(3, 3).Value = 4
$totalSheet.Cells.Item(3, 4).Value = 0.05

$totalSheet.Cells.Item(4, 1).Value = 2
